# Apply the "updated lat lon outputs for TT rule + updated scripts for TRI
# facilities that dont drop facilities" changes to the facility demographics
# worksheet.
#
# Net effect (validated against the target OOXML diff):
#   - Row 2   (Linde-Decatur / Decatur)            -> refreshed metrics
#   - Row 6/7 swap identities: the facility that used to be row 6
#     (Linde-Whiting / EAST CHICAGO) is now row 7 (Linde-Whiting / East
#     Chicago - re-cased), and the facility that used to be row 7
#     (Diversified-CPC / Channahon) is now row 6, each row keeping its own
#     refreshed metrics and the GHG_co2e (col C) value travelling with the
#     Linde-Whiting facility (row 6 -> row 7).
#   - City-name casing normalized: SERGEANT BLUFF -> Sergeant Bluff,
#     GEISMAR -> Geismar (rows 8 and 9 respectively).
#   - Rows 11, 13, 14 get refreshed metrics (labels/cities unchanged).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 2: Linde-Decatur / Decatur ----
$ws.Range("D2").Value = 18
$ws.Range("E2").Value = 171.012007267968
$ws.Range("F2").Value = 26593
$ws.Range("G2").Value = 155.503700733306
$ws.Range("I2").Value = 0.17
$ws.Range("J2").Value = 15888
$ws.Range("K2").Value = 8563
$ws.Range("L2").Value = 211
$ws.Range("M2").Value = 60
$ws.Range("N2").Value = 3609
$ws.Range("O2").Value = 46.6776666666667
$ws.Range("P2").Value = 8.49845991335092
$ws.Range("Q2").Value = 14.2254707164624
$ws.Range("R2").Value = 40.5555555555556
$ws.Range("S2").Value = 0.505555555555556

# ---- Row 6: becomes Diversified-CPC / Channahon ----
$ws.Range("A6").Value = "Diversified-CPC"
$ws.Range("B6").Value = "Channahon"
$ws.Range("C6").Value = $null
$ws.Range("D6").Value = 11
$ws.Range("E6").Value = 132.596145486872
$ws.Range("F6").Value = 24316
$ws.Range("G6").Value = 183.383912938913
$ws.Range("I6").Value = 0.36
$ws.Range("J6").Value = 22281
$ws.Range("K6").Value = 493
$ws.Range("L6").Value = 42
$ws.Range("M6").Value = 199
$ws.Range("N6").Value = 3206
$ws.Range("O6").Value = 94.4638
$ws.Range("P6").Value = 1.96299114127479
$ws.Range("Q6").Value = 2.94731126966235
$ws.Range("R6").Value = 27
$ws.Range("S6").Value = 0.32

# ---- Row 7: becomes Linde-Whiting / East Chicago ----
$ws.Range("A7").Value = "Linde-Whiting"
$ws.Range("B7").Value = "East Chicago"
$ws.Range("C7").Value = 183562
$ws.Range("D7").Value = 65
$ws.Range("E7").Value = 39.1884589201756
$ws.Range("F7").Value = 70220
$ws.Range("G7").Value = 1791.85407987167
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 24499
$ws.Range("K7").Value = 20429
$ws.Range("L7").Value = 176
$ws.Range("M7").Value = 490
$ws.Range("N7").Value = 34175
$ws.Range("O7").Value = 39.5717419354839
$ws.Range("P7").Value = 12.9769263034948
$ws.Range("Q7").Value = 13.897701347952
$ws.Range("R7").Value = 30
$ws.Range("S7").Value = 0.364615384615385

# ---- Row 8: CFI-PortNeal / Sergeant Bluff (re-cased city) ----
$ws.Range("B8").Value = "Sergeant Bluff"

# ---- Row 9: APC-Geismar / Geismar (re-cased city) ----
$ws.Range("B9").Value = "Geismar"

# ---- Row 11: AEROPRES-SIBLEY / Sibley ----
$ws.Range("D11").Value = 8
$ws.Range("E11").Value = 155.317267451068
$ws.Range("F11").Value = 9240
$ws.Range("G11").Value = 59.4911316149121
$ws.Range("H11").Value = 0
$ws.Range("I11").Value = 0.75
$ws.Range("J11").Value = 4737
$ws.Range("K11").Value = 4379
$ws.Range("L11").Value = 7
$ws.Range("M11").Value = 37
$ws.Range("N11").Value = 162
$ws.Range("O11").Value = 25.380125
$ws.Range("P11").Value = 9.17561925571522
$ws.Range("Q11").Value = 21.1990327607574

# ---- Row 13: APC-PortAuthur / Port Arthur ----
$ws.Range("D13").Value = 19
$ws.Range("E13").Value = 449.535365707168
$ws.Range("F13").Value = 14855
$ws.Range("G13").Value = 33.0452309945213
$ws.Range("I13").Value = 0.32
$ws.Range("J13").Value = 4680
$ws.Range("K13").Value = 9373
$ws.Range("L13").Value = 43
$ws.Range("M13").Value = 237
$ws.Range("N13").Value = 2636
$ws.Range("O13").Value = 34.5933888888889
$ws.Range("P13").Value = 16.174443514373
$ws.Range("Q13").Value = 14.4631067590132
$ws.Range("R13").Value = 42.1052631578947

# ---- Row 14: Chemours-CorpusChristie / Gregory ----
$ws.Range("D14").Value = 11
$ws.Range("E14").Value = 322.88741866344
$ws.Range("F14").Value = 19446
$ws.Range("G14").Value = 60.2253258442052
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0.91
$ws.Range("J14").Value = 17821
$ws.Range("K14").Value = 420
$ws.Range("L14").Value = 86
$ws.Range("M14").Value = 167
$ws.Range("N14").Value = 8387
$ws.Range("O14").Value = 73.1127
$ws.Range("P14").Value = 4.34232462474073
$ws.Range("Q14").Value = 3.89252789332616
$ws.Range("S14").Value = 0.209090909090909

Write-Output "Edit applied."
